# Auto-generated: apply 2022-05-06 crime data updates to column I (year 2022) values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 2085
$ws.Range("I3").Value = 2196
$ws.Range("I4").Value = 552
$ws.Range("I5").Value = 195
$ws.Range("I6").Value = 2628
$ws.Range("I7").Value = 7656

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 81
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 260
$ws.Range("I8").Value = 480
$ws.Range("I13").Value = 11
$ws.Range("I15").Value = 97
$ws.Range("I17").Value = 6
$ws.Range("I19").Value = 215
$ws.Range("I20").Value = 206
$ws.Range("I23").Value = 64
$ws.Range("I29").Value = 491
$ws.Range("I33").Value = 357
$ws.Range("I36").Value = 100
$ws.Range("I37").Value = 245
$ws.Range("I42").Value = 257
$ws.Range("I44").Value = 61
$ws.Range("I48").Value = 79
$ws.Range("I49").Value = 46
$ws.Range("I52").Value = 154
$ws.Range("I53").Value = 77
$ws.Range("I54").Value = 181
$ws.Range("I55").Value = 87
$ws.Range("I57").Value = 23
$ws.Range("I63").Value = 30
$ws.Range("I67").Value = 298
$ws.Range("I76").Value = 119
$ws.Range("I79").Value = 198
$ws.Range("I82").Value = 7
$ws.Range("I84").Value = 54
$ws.Range("I85").Value = 358
$ws.Range("I88").Value = 62
$ws.Range("I89").Value = 80
$ws.Range("I90").Value = 88
$ws.Range("I92").Value = 23
$ws.Range("I94").Value = 64
$ws.Range("I96").Value = 100
$ws.Range("I97").Value = 60
$ws.Range("I100").Value = 10
$ws.Range("I101").Value = 7656

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 88
$ws.Range("I3").Value = 138
$ws.Range("I7").Value = 358

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I3").Value = 62
$ws.Range("I6").Value = 31
$ws.Range("I7").Value = 154

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I3").Value = 130
$ws.Range("I4").Value = 28
$ws.Range("I6").Value = 159
$ws.Range("I7").Value = 480

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I3").Value = 23
$ws.Range("I7").Value = 77

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 76
$ws.Range("I7").Value = 260

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I2").Value = 22
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 80

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I2").Value = 25
$ws.Range("I7").Value = 100

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I4").Value = 20
$ws.Range("I6").Value = 69
$ws.Range("I7").Value = 245

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 105
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 298

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I3").Value = 18
$ws.Range("I7").Value = 54

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 86
$ws.Range("I6").Value = 121
$ws.Range("I7").Value = 357

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 46

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I4").Value = 15
$ws.Range("I7").Value = 181

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 151
$ws.Range("I3").Value = 162
$ws.Range("I6").Value = 149
$ws.Range("I7").Value = 491

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 91
$ws.Range("I7").Value = 215

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I2").Value = 21
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I3").Value = 29
$ws.Range("I7").Value = 119

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 21
$ws.Range("I7").Value = 50

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 68
$ws.Range("I7").Value = 257

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 11

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value = 28
$ws.Range("I3").Value = 23
$ws.Range("I7").Value = 87

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 56
$ws.Range("I7").Value = 198

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 53
$ws.Range("I3").Value = 62
$ws.Range("I7").Value = 206

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("I4").Value = 2
$ws.Range("I7").Value = 6

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I3").Value = 33
$ws.Range("I7").Value = 100

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 10

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 28
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 24
$ws.Range("I5").Value = 2
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 81

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I3").Value = 9
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 60

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I4").Value = 2
$ws.Range("I7").Value = 23

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 88

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 23

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("I5").Value = 3
$ws.Range("I6").Value = 7
